$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19, shifting existing rows (19-39) down to (20-40)
$ws.Rows("19:19").Insert()

# Populate the new row 19 with the NWWv4 (ammonia water world) data
$ws.Range("B19").Value = "NWWv4"
$ws.Range("C19").Value = "WW"
$ws.Range("D19").Value = 387
$ws.Range("E19").Value = "none"
$ws.Range("F19").Value = "ammonia"
$ws.Range("G19").Value = "ammonia"
$ws.Range("H19").Value = "nitrogen"
$ws.Range("I19").Value = "oxigen"
$ws.Range("J19").Value = 67.3
$ws.Range("K19").Value = 32.700000000000003

# Match the row height used by the other 18.75pt data rows in this table
$ws.Rows("19:19").RowHeight = 18.75

# Update the selected cell to match the new active selection
$ws.Range("K19").Select()
